# Deduplicate rows in the key_advantages workbook.
#
# The source data for several banks contained duplicate advantage rows
# (same advantage_id/title/category/evidence_list, differing only in the
# example_ad_index / example_quote columns). This edit removes the
# redundant rows so each advantage keeps only the intended example row(s):
#   - Artea:               drop original rows 3 and 5
#   - Luminor_Lietuva:     drop original rows 3, 4, 7 and 9
#   - Swedbank_Lietuvoje:  drop original row 5
# Citadele_bankas and SEB_Lietuvoje are untouched.
#
# Rows are deleted bottom-up within each sheet so earlier deletions don't
# shift the row numbers of rows still queued for removal.

$wb = $excel.ActiveWorkbook

# --- Artea ---------------------------------------------------------------
$wsArtea = $wb.Worksheets.Item("Artea")
$wsArtea.Rows.Item(5).Delete()
$wsArtea.Rows.Item(3).Delete()

# --- Luminor_Lietuva -------------------------------------------------------
$wsLuminor = $wb.Worksheets.Item("Luminor_Lietuva")
$wsLuminor.Rows.Item(9).Delete()
$wsLuminor.Rows.Item(7).Delete()
$wsLuminor.Rows.Item(4).Delete()
$wsLuminor.Rows.Item(3).Delete()

# --- Swedbank_Lietuvoje ----------------------------------------------------
$wsSwedbank = $wb.Worksheets.Item("Swedbank_Lietuvoje")
$wsSwedbank.Rows.Item(5).Delete()

# --- Selections / active sheet, matching the final saved view -------------
[void]$wsArtea.Range("E12").Select()

$wsLuminor.Range("G13").Select() | Out-Null

$wsSwedbank.Activate()
$wsSwedbank.Range("F8").Select() | Out-Null
